$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column J (old J "TOTAL" shifts right to L),
# carrying along styles/merges/row spans the way Excel's own column-insert does.
$ws.Columns("J:K").Insert() | Out-Null

# New retention headers in the freshly inserted columns.
$ws.Range("J9").Value = "RENTA"
$ws.Range("K9").Value = "IVA"

# Match the header formatting used by the other column headers in row 9.
$ws.Range("B9").Copy() | Out-Null
$ws.Range("J9:K9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Give the two new columns the same width as column I.
$ws.Columns.Item(10).ColumnWidth = $ws.Columns.Item(9).ColumnWidth
$ws.Columns.Item(11).ColumnWidth = $ws.Columns.Item(9).ColumnWidth

# Move the active selection to the new RENTA header cell.
$ws.Range("J9").Select() | Out-Null
